$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the street name spelling: "A.J. Ernststraat" -> "A. J. Ernststraat"
$ws.Range("A8").Value = "A. J. Ernststraat"

# Move the active selection to A9
[void]$ws.Range("A9").Select()
